$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.720121383666992
$ws.Range("B1").Value = 5.971473217010498
$ws.Range("C1").Value = 6.307316780090332
$ws.Range("D1").Value = 9.894746780395508
$ws.Range("E1").Value = 7.229799270629883
